$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find a paragraph whose text starts with $needle and return its
# Range (paragraph mark excluded) so InsertXML only touches the run content
# and leaves <w:pPr> (paragraph formatting) untouched.
# ---------------------------------------------------------------------------
function Get-ParaContentRange($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($needle)) {
            $full = $p.Range
            return $d.Range($full.Start, $full.End - 1)
        }
    }
    return $null
}

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$xmlFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# Paragraph: "STUDENT NAME:"  ->  "STUDENT NAME: Sanaullah"
# Word split "STUDENT NAME" into "STUDENT " + "NAME" (spell-check boundary)
# and wrapped the inserted answer ": Sanaullah" with grammar-check markers.
# ---------------------------------------------------------------------------
$nameRange = Get-ParaContentRange("STUDENT NAME")
if ($null -eq $nameRange) { throw "Could not locate the 'STUDENT NAME' paragraph" }
$nameBody = @'
<w:body><w:p>
<w:r w:rsidRPr="00B46F38"><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">STUDENT </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>NAME</w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r>
<w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:sz w:val="44"/><w:szCs w:val="28"/></w:rPr><w:t>Sanaullah</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:proofErr w:type="gramEnd"/>
</w:p></w:body>
'@
$nameRange.InsertXML($xmlHeader + $nameBody + $xmlFooter)

# ---------------------------------------------------------------------------
# Paragraph: "STUDENT ID:"  ->  "STUDENT ID: 22066704"
# Word split "STUDENT ID:" into "STUDENT ID" + ":" and wrapped the inserted
# answer ": 22066704" with grammar-check markers.
# ---------------------------------------------------------------------------
$idRange = Get-ParaContentRange("STUDENT ID")
if ($null -eq $idRange) { throw "Could not locate the 'STUDENT ID' paragraph" }
$idBody = @'
<w:body><w:p>
<w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>STUDENT ID</w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r>
<w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="44"/><w:szCs w:val="28"/></w:rPr><w:t>22066704</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
</w:p></w:body>
'@
$idRange.InsertXML($xmlHeader + $idBody + $xmlFooter)

Write-Output "done"
